$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-05-27 Monday" "2024-05-28 Tuesday"

Replace-Text "25÷3=8, 1" "78÷9=8, 6"
Replace-Text "63÷6=10, 3" "98÷3=32, 2"
Replace-Text "64÷5=12, 4" "24÷3=8, 0"
Replace-Text "16÷9=1, 7" "18÷3=6, 0"
Replace-Text "87÷4=21, 3" "36÷3=12, 0"
Replace-Text "55÷5=11, 0" "57÷6=9, 3"
Replace-Text "92÷2=46, 0" "55÷6=9, 1"
Replace-Text "97÷7=13, 6" "31÷4=7, 3"
Replace-Text "64÷9=7, 1" "23÷9=2, 5"
Replace-Text "41÷3=13, 2" "86÷8=10, 6"
Replace-Text "94÷6=15, 4" "12÷8=1, 4"
Replace-Text "91÷2=45, 1" "67÷7=9, 4"
Replace-Text "34÷4=8, 2" "46÷3=15, 1"
Replace-Text "89÷9=9, 8" "27÷6=4, 3"
Replace-Text "27÷9=3, 0" "33÷7=4, 5"
Replace-Text "65÷8=8, 1" "38÷2=19, 0"
Replace-Text "71÷9=7, 8" "67÷9=7, 4"
Replace-Text "12÷6=2, 0" "66÷4=16, 2"
Replace-Text "70÷2=35, 0" "29÷5=5, 4"
Replace-Text "93÷5=18, 3" "42÷4=10, 2"
Replace-Text "74÷5=14, 4" "66÷2=33, 0"
Replace-Text "48÷3=16, 0" "99÷4=24, 3"
Replace-Text "38÷4=9, 2" "46÷4=11, 2"
Replace-Text "85÷2=42, 1" "18÷3=6, 0"
Replace-Text "35÷9=3, 8" "15÷3=5, 0"

Write-Host "Done applying replacements"
